# Invitation task MCQ mapping workbook update
# - Re-shuffles several answer-option cells (columns E-H) across the trial
#   table on Sheet1 so the stimulus mapping matches the new schedule.
# - Two new distractor values ("20-20" and "10-25") are introduced, which
#   Excel will automatically append to the shared-strings table.
# - Updates the active cell selection left over from the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "5-20"
$ws.Range("G3").Value = "0-5"
$ws.Range("G4").Value = "0-5"
$ws.Range("H4").Value = "0-20"
$ws.Range("E5").Value = "5-20"
$ws.Range("G7").Value = "0-0"
$ws.Range("F8").Value = "20-20"
$ws.Range("E9").Value = "5-10"
$ws.Range("F10").Value = "20-20"
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "5-2"
$ws.Range("E12").Value = "10-0"
$ws.Range("G12").Value = "5-5"
$ws.Range("F13").Value = "0-0"
$ws.Range("H14").Value = "5-10"
$ws.Range("F15").Value = "10-5"
$ws.Range("G15").Value = "0-5"
$ws.Range("G16").Value = "0-0"
$ws.Range("E17").Value = "10-20"
$ws.Range("G19").Value = "5-10"
$ws.Range("E20").Value = "10-25"
$ws.Range("G21").Value = "5-20"
$ws.Range("E22").Value = "10-20"
$ws.Range("E23").Value = "0-0"
$ws.Range("H24").Value = "10-5"
$ws.Range("E25").Value = "5-10"
$ws.Range("G26").Value = "5-20"
$ws.Range("H27").Value = "10-25"
$ws.Range("E28").Value = "20-0"
$ws.Range("F28").Value = "5-10"
$ws.Range("E29").Value = "20-0"

# Move the saved selection to match where the author left off editing.
$ws.Range("J16").Select()
